# Updated cryptos list on Mon May  6 17:39:47 UTC 2024 with GitHub Actions
# Refreshes Price (col D) / Volume(1h) (col E) for every coin row, and for
# a few coins whose rank swapped with their neighbor, also updates the
# Coin name (col B) and Link (col C) so each row keeps matching data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.273.44"
$ws.Range("E2").Value = "  -1.98%  "

# Row 3
$ws.Range("D3").Value = "3.080.86"
$ws.Range("E3").Value = "  -2.69%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'587.20"
$ws.Range("E5").Value = "  -1.00%  "

# Row 6
$ws.Range("D6").Value = "'151.92"
$ws.Range("E6").Value = "  +3.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "'0.555"
$ws.Range("E8").Value = "  +4.44%  "

# Row 9
$ws.Range("D9").Value = "3.069.51"
$ws.Range("E9").Value = "  -2.72%  "

# Row 10
$ws.Range("D10").Value = "'0.155"
$ws.Range("E10").Value = "  -4.42%  "

# Row 11
$ws.Range("D11").Value = "'5.85"
$ws.Range("E11").Value = "  -1.31%  "

# Row 12
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  -0.43%  "

# Row 13
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "  -3.33%  "

# Row 14
$ws.Range("D14").Value = "'37.09"
$ws.Range("E14").Value = "  -0.95%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.591.01"
$ws.Range("E15").Value = "  -2.69%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.119"
$ws.Range("E16").Value = "  -2.01%  "

# Row 17
$ws.Range("D17").Value = "'7.18"
$ws.Range("E17").Value = "  -1.08%  "

# Row 18
$ws.Range("D18").Value = "63.331.04"
$ws.Range("E18").Value = "  -1.53%  "

# Row 19
$ws.Range("D19").Value = "3.076.66"
$ws.Range("E19").Value = "  -2.59%  "

# Row 20
$ws.Range("D20").Value = "'479.19"
$ws.Range("E20").Value = "  +2.00%  "

# Row 21
$ws.Range("D21").Value = "'14.62"
$ws.Range("E21").Value = "  +0.82%  "

# Row 22
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  -2.80%  "

# Row 23
$ws.Range("D23").Value = "'7.51"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24
$ws.Range("D24").Value = "'13.03"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  -0.86%  "

# Row 26
$ws.Range("D26").Value = "'81.50"
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.24%  "

# Row 28
$ws.Range("D28").Value = "'9.66"
$ws.Range("E28").Value = "  -1.19%  "

# Row 29
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.67"
$ws.Range("E30").Value = "  -2.49%  "

# Row 31
$ws.Range("D31").Value = "'2.19"
$ws.Range("E31").Value = "  -3.75%  "

# Row 32
$ws.Range("D32").Value = "'7.18"
$ws.Range("E32").Value = "  -3.43%  "

# Row 33
$ws.Range("E33").Value = "  +1.51%  "

# Row 34
$ws.Range("D34").Value = "'27.28"
$ws.Range("E34").Value = "  -1.55%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("E36").Value = "  -2.36%  "

# Row 37
$ws.Range("D37").Value = "'6.10"
$ws.Range("E37").Value = "  -3.04%  "

# Row 38
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  +2.43%  "

# Row 39
$ws.Range("D39").Value = "'2.20"
$ws.Range("E39").Value = "  -5.51%  "

# Row 40
$ws.Range("D40").Value = "'9.35"
$ws.Range("E40").Value = "  +0.52%  "

# Row 41
$ws.Range("D41").Value = "'50.46"
$ws.Range("E41").Value = "  -2.47%  "

# Row 42
$ws.Range("D42").Value = "'440.29"
$ws.Range("E42").Value = "  -7.24%  "

# Row 43
$ws.Range("D43").Value = "'0.283"
$ws.Range("E43").Value = "  -4.30%  "

# Row 44
$ws.Range("D44").Value = "'0.0363"
$ws.Range("E44").Value = "  -3.51%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.807.53"
$ws.Range("E45").Value = "  -4.32%  "

# Row 46
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "'39.67"
$ws.Range("E46").Value = "  -2.84%  "

# Row 47
$ws.Range("E47").Value = "  +0.60%  "

# Row 48
$ws.Range("D48").Value = "'130.57"
$ws.Range("E48").Value = "  +0.75%  "

# Row 50
$ws.Range("D50").Value = "'0.112"
$ws.Range("E50").Value = "  +0.90%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.25"
$ws.Range("E51").Value = "  -0.78%  "
